$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adduser")

$ws.Range("A5").Value = "ertrty"
$ws.Range("B5").Value = 23435355

$ws.Range("A6").Value = "keva"
$ws.Range("B6").Value = 132154667

$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

$ws.Range("B6").Select() | Out-Null
